# Mapping-file-creator template update:
# "primers with 3 or 4 digits" - expand the FILL_IN example table from a
# handful of placeholder FW00x/RV00x primer names to the full set of real
# 515F / 926RBC Golay-barcoded primer names, paired up or down with sample
# names S1..S33 in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FILL_IN")

$forward = @("515F_Golay0241","515F_Golay0243","515F_Golay0245","515F_Golay0247","515F_Golay0249","515F_Golay0251","515F_Golay0253","515F_Golay0255","515F_Golay0401","515F_Golay0402","515F_Golay0403","515F_Golay0404","515F_Golay9001","515F_Golay9003","515F_Golay9005","515F_Golay9007","515F_Golay9009","515F_Golay9011","515F_Golay9013","515F_Golay9015","515F_Golay9017","515F_Golay9019","515F_Golay9021","515F_Golay9023","515F_Golay9025","515F_Golay9027","515F_Golay9029","515F_Golay9031","515F_Golay9033","515F_Golay9035","515F_Golay9037","515F_Golay9039","515F_Golay9041")

$reverse = @("926RBC_Golay0160","926RBC_Golay0162","926RBC_Golay0164","926RBC_Golay0166","926RBC_Golay0168","926RBC_Golay0170","926RBC_Golay9046","926RBC_Golay9048","926RBC_Golay9050","926RBC_Golay9052","926RBC_Golay9054","926RBC_Golay9056","926RBC_Golay9058","926RBC_Golay9060","926RBC_Golay9062","926RBC_Golay9064","926RBC_Golay0024","926RBC_Golay0026","926RBC_Golay0028","926RBC_Golay0030","926RBC_Golay0032","926RBC_Golay0034","926RBC_Golay0036","926RBC_Golay9106","926RBC_Golay9108","926RBC_Golay9110","926RBC_Golay9112","926RBC_Golay9114","926RBC_Golay9116","926RBC_Golay9118","926RBC_Golay9120","926RBC_Golay9122","926RBC_Golay9124")

$samples = @("S1","S2","S3","S4","S5","S6","S7","S8","S9","S10","S11","S12","S13","S14","S15","S16","S17","S18","S19","S20","S21","S22","S23","S24","S25","S26","S27","S28","S29","S30","S31","S32","S33")

# Rows 2..34 hold one forward/reverse primer pair (columns A/B) plus one
# sample name (column D). Columns A/B go back to the sheet's default
# (unstyled) look, column D keeps the data-entry ("unlocked") style that
# the column already carries.
for ($i = 0; $i -lt $forward.Length; $i++) {
    $row = 2 + $i

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value2 = $forward[$i]
    $cellA.Style = "Normal"

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value2 = $reverse[$i]
    $cellB.Style = "Normal"

    $ws.Cells.Item($row, 4).Value2 = $samples[$i]
}

# Trailing blank spacer rows below the table (35-38), present in the sheet
# but carrying no value/style of their own.
for ($row = 35; $row -le 38; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value2 = ""
    $cell.Style = "Normal"
}

$ws.Activate()
$ws.Range("J23").Select()
